# Acme Chollos y Rifas - minor requirement wording fixes
# (Modificaciones menores en los requisitos.)

$d = $word.ActiveDocument

# 1) Fix typo "des"/"een." split in the wishlist bullet (merges the two runs
#    and drops the stray _GoBack bookmark Word had left sitting mid-word).
$d.Content.Find.Execute(
    "que deseen.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "que deseen.", 2) | Out-Null

# 2) Survey-restriction paragraph: tighten the actor wording, reword the
#    "moderator filters users" sentence to a plain statement, and add a new
#    sentence about companies opting sponsors with ads into surveys.
$d.Content.Find.Execute(
    "(patrocinador, usuario o ambos)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(patrocinador o usuario)", 2) | Out-Null

$d.Content.Find.Execute(
    "Si el autor es un moderador puede filtrar los usuarios",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Los usuarios pueden ser filtrados", 2) | Out-Null

$d.Content.Find.Execute(
    "Una vez publicada la encuesta",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Las compañías pueden definir si quieres que sus encuestas lleguen a los patrocinadores que tienen anuncios en sus chollos. Una vez publicada la encuesta",
    2) | Out-Null

$d.Content.Find.Execute(
    "los actores seleccionados según las restricciones deben ser notificados",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "los actores seleccionados deben ser notificados", 2) | Out-Null

# 3) Points paragraph: participating in a "chollo" no longer grants points,
#    only joining a "conjunta" does.
$d.Content.Find.Execute(
    "participa en un chollo o conjunta",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "participa en una conjunta", 2) | Out-Null
